$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(3)

# --- Row 7 gets the OLD "Movetype" field (moved down from row 6) ---
# Write the label first, while the original AG6 still references the shared
# string "Movetype", so it is reused instead of creating a duplicate entry.
$ws.Range("D7").Value = 1
$ws.Range("E7").Value = 1
$ws.Range("F7").Value = 1
$ws.Range("AG7").Value = "Movetype"
$ws.Range("AH7").Value = 7
$ws.Range("AI7").Formula = "=HEX2BIN(AH7)"
$ws.Range("AJ7").Formula = "=HEX2DEC(AH7)"

# --- Row 6 becomes the new "Promotion" field ---
$ws.Range("G6").Value = 1
$ws.Range("AG6").Value = "Promotion"
$ws.Range("AH6").Value = "f"

# --- Row 8: new "Castling" field ---
$ws.Range("A8").Value = 1
$ws.Range("B8").Value = 1
$ws.Range("C8").Value = 1
$ws.Range("AG8").Value = "Castling"
$ws.Range("AH8").Value = 7
$ws.Range("AI8").Formula = "=HEX2BIN(AH8)"
$ws.Range("AJ8").Formula = "=HEX2DEC(AH8)"

# --- Row 9: full-byte ("ff") sanity check row ---
$ws.Range("AH9").Value = "ff"
$ws.Range("AI9").Formula = "=HEX2BIN(AH9)"
$ws.Range("AJ9").Formula = "=HEX2DEC(AH9)"

# Column AI now holds wider text ("11111111"); resize it to fit the content.
$ws.Range("AI1").EntireColumn.AutoFit() | Out-Null

# Move the active selection, matching where editing ended up.
$ws.Range("O16").Select() | Out-Null
